# Generate Report for Handback
#
# Updates the localization-status report after a handback event:
#  - Status moves from "Ready for handoff" -> "Handed back: in sync with en-US"
#    (this text is a shared string used by Overview!E2/F2 and the per-locale
#    sheets' Status cell, so it is set on every cell that shows it)
#  - Each locale sheet's "Latest Target File" (I2) and "Latest Handback File"
#    (J2) columns get populated, with I2 becoming a hyperlink to the source
#    markdown file (mirroring the existing A2 hyperlink)
#  - Each locale sheet's "Latest Handback DateTime" (K2) is stamped with the
#    actual handback time
#  - A few columns are widened to fit the newly-populated long filenames

$wb = $excel.ActiveWorkbook

$ws_overview = $wb.Worksheets.Item("Overview")
$ws_zhcn = $wb.Worksheets.Item("zh-cn")
$ws_dede = $wb.Worksheets.Item("de-de")

$statusText = "Handed back: in sync with en-US"

# --- Status text: flips for both locales everywhere it is shown ---
$ws_overview.Range("E2").Value = $statusText
$ws_overview.Range("F2").Value = $statusText
$ws_zhcn.Range("C2").Value = $statusText
$ws_dede.Range("C2").Value = $statusText

# --- zh-cn locale sheet: target/handback file + datetime ---
$ws_zhcn.Range("I2").Value = "759239cf-2818-4f34-9c84-0fae4df38b1c.md"
$ws_zhcn.Range("I2").Style = "Hyperlink"
$ws_zhcn.Hyperlinks.Add(
    $ws_zhcn.Range("I2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2978d00b5a530dcd104d2590d7e7202a99367c54/e2e/759239cf-2818-4f34-9c84-0fae4df38b1c.md",
    "",
    "",
    "759239cf-2818-4f34-9c84-0fae4df38b1c.md"
) | Out-Null
$ws_zhcn.Range("J2").Value = "759239cf-2818-4f34-9c84-0fae4df38b1c.229120e0fd521357ec57b6b626152d250f46545a.zh-cn.xlf"
$ws_zhcn.Range("K2").Value = "2016-08-19 15:05:23"

# --- de-de locale sheet: target/handback file + datetime ---
$ws_dede.Range("I2").Value = "759239cf-2818-4f34-9c84-0fae4df38b1c.md"
$ws_dede.Range("I2").Style = "Hyperlink"
$ws_dede.Hyperlinks.Add(
    $ws_dede.Range("I2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2978d00b5a530dcd104d2590d7e7202a99367c54/e2e/759239cf-2818-4f34-9c84-0fae4df38b1c.md",
    "",
    "",
    "759239cf-2818-4f34-9c84-0fae4df38b1c.md"
) | Out-Null
$ws_dede.Range("J2").Value = "759239cf-2818-4f34-9c84-0fae4df38b1c.229120e0fd521357ec57b6b626152d250f46545a.de-de.xlf"
$ws_dede.Range("K2").Value = "2016-08-19 15:05:30"

# --- Column widths: widen columns holding the newly-populated long file
#     names/paths. The engine snaps ColumnWidth to a whole-pixel grid, so
#     feed it (desired - 5/6) to land on the nearest achievable width. ---
$wide = 30 - (5 / 6)
$wider = 40 - (5 / 6)

$ws_overview.Range("E1").ColumnWidth = $wide
$ws_overview.Range("F1").ColumnWidth = $wide

$ws_zhcn.Range("C1").ColumnWidth = $wide
$ws_zhcn.Range("I1").ColumnWidth = $wider
$ws_zhcn.Range("J1").ColumnWidth = $wider

$ws_dede.Range("C1").ColumnWidth = $wide
$ws_dede.Range("I1").ColumnWidth = $wider
$ws_dede.Range("J1").ColumnWidth = $wider
